$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.834
$ws.Range("C3").Value = -12.093
$ws.Range("E4").Value = 12.662
$ws.Range("C5").Value = -12.634
$ws.Range("E6").Value = 12.338
$ws.Range("D7").Value = -7.237
$ws.Range("B9").Value = 6.658999999999999
$ws.Range("D9").Value = -8.109999999999999
$ws.Range("E10").Value = 12.38
$ws.Range("C11").Value = -12.628
$ws.Range("E11").Value = 13.072
$ws.Range("C12").Value = -12.628
$ws.Range("B13").Value = 6.307999999999999
$ws.Range("B16").Value = 5.786999999999999
$ws.Range("B18").Value = 5.91
$ws.Range("B20").Value = 6.308000000000001
$ws.Range("C21").Value = -12.09
$ws.Range("D21").Value = -7.734999999999999
$ws.Range("E21").Value = 12.961
$ws.Range("E25").Value = 12.763
